# Updated update_count function behaviour:
#   - "add" actions append a new activity row to the "<Section> Timestamps"
#     log sheet and push the SAN into the "All SANs" lookup sheet.
#   - "subtract" actions now check whether the SAN exists in "All SANs"; if it
#     does, that SAN's row is removed from "All SANs" (the rows below shift
#     up) in addition to being logged in the timestamps sheet.
#
# This script replays that batch of activity (5 adds + 1 subtract of
# "SAN111111", all for "Desktop Mini G9") against the "4.2 *" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "4.2 Items" - refreshed LastCount/NewCount for the affected item rows.
# ---------------------------------------------------------------------------
$itemsWs = $wb.Worksheets.Item("4.2 Items")

$itemsWs.Cells.Item(2, 2).Value = 110   # Desktop Mini G9 - LastCount
$itemsWs.Cells.Item(2, 3).Value = 109   # Desktop Mini G9 - NewCount

$itemsWs.Cells.Item(4, 2).Value = 5     # Dock Thunderbolt G2 - LastCount
$itemsWs.Cells.Item(4, 3).Value = 7     # Dock Thunderbolt G2 - NewCount

$itemsWs.Cells.Item(7, 3).Value = 34    # Laptop 840 G9 - NewCount

# ---------------------------------------------------------------------------
# 2) "4.2 Timestamps" - append the new activity log rows (28-33).
#    Rows 28-32 sit inside the already-formatted table, so clone the last
#    existing data row (27) to inherit its per-cell formatting before
#    overwriting the values. Row 33 is appended past the end completely
#    unformatted (matches how the source tool wrote it).
# ---------------------------------------------------------------------------
$tsWs = $wb.Worksheets.Item("4.2 Timestamps")

$addedRows = @(
  @("2023-12-31 12:35:48", "Desktop Mini G9", "add", "SAN122334"),
  @("2023-12-31 12:36:48", "Desktop Mini G9", "add", "SAN222211"),
  @("2023-12-31 12:41:56", "Desktop Mini G9", "add", "SAN456789"),
  @("2023-12-31 12:41:56", "Desktop Mini G9", "add", "SAN125689"),
  @("2023-12-31 12:41:57", "Desktop Mini G9", "add", "SAN357895")
)

$templateRow = 27
for ($i = 0; $i -lt $addedRows.Count; $i++) {
  $destRow = $templateRow + 1 + $i
  $tsWs.Rows.Item($templateRow).Copy()
  $tsWs.Rows.Item($destRow).Insert()
}

for ($i = 0; $i -lt $addedRows.Count; $i++) {
  $r = $templateRow + 1 + $i
  $rowData = $addedRows[$i]
  $tsWs.Cells.Item($r, 1).Value = $rowData[0]
  $tsWs.Cells.Item($r, 2).Value = $rowData[1]
  $tsWs.Cells.Item($r, 3).Value = $rowData[2]
  $tsWs.Cells.Item($r, 4).Value = $rowData[3]
}

# Final "subtract" row - appended with no direct formatting at all, so stage
# it on a throwaway sheet (no column styles defined there) and Cut/Paste it
# into place; that preserves the "no style" cells instead of the target
# columns' inherited default column format.
$subtractRow = 33
$scratch = $wb.Worksheets.Add()
$scratch.Cells.Item(1, 1).Value = "2023-12-31 14:20:35"
$scratch.Cells.Item(1, 2).Value = "Desktop Mini G9"
$scratch.Cells.Item(1, 3).Value = "subtract"
$scratch.Cells.Item(1, 4).Value = "SAN111111"
$scratch.Range("A1:D1").Cut($tsWs.Range("A" + $subtractRow + ":D" + $subtractRow))
$excel.DisplayAlerts = $false
[void]$scratch.Delete()

# ---------------------------------------------------------------------------
# 3) "All SANs" - the "subtract" of SAN111111 removes that entry. Row 8 held
#    SAN111111, so it (and the rows after it) are replaced/extended with the
#    still-outstanding "add" entries, in order.
# ---------------------------------------------------------------------------
$sanWs = $wb.Worksheets.Item("All SANs")

$pendingAdds = @(
  @("Desktop Mini G9", "SAN122334", "2023-12-31 12:35:48"),
  @("Desktop Mini G9", "SAN222211", "2023-12-31 12:36:48"),
  @("Desktop Mini G9", "SAN456789", "2023-12-31 12:41:56"),
  @("Desktop Mini G9", "SAN125689", "2023-12-31 12:41:56"),
  @("Desktop Mini G9", "SAN357895", "2023-12-31 12:41:57")
)

$sanStartRow = 8
for ($i = 0; $i -lt $pendingAdds.Count; $i++) {
  $r = $sanStartRow + $i
  $rowData = $pendingAdds[$i]
  $sanWs.Cells.Item($r, 1).Value = $rowData[0]
  $sanWs.Cells.Item($r, 2).Value = $rowData[1]
  $sanWs.Cells.Item($r, 3).Value = $rowData[2]
}

Write-Host "Edit complete"
